# "data model of structures"
# Rework the "structures" sheet header/columns to the new INSPER_* data
# model (20 columns, A:T), and rename the "structure_name" header on the
# "sections" sheet to "BUSINESS_TITLE".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "structures" sheet -> new INSPER_* column layout
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("structures")

# Stretch the existing header formatting (bold font, thin border, centered)
# across the new D1:T1 header cells before touching any values, by copying
# the format of the already-styled A1:C1 range.
$ws.Range("A1:C1").Copy()
$ws.Range("D1:T1").PasteSpecial(-4122)

$headers = @(
    "INSPER_ID_PRE",
    "BUSINESS_ID_PRE",
    "TYPE_OF_PARTICIPATION_CD",
    "TYPE_OF_INSURED_PERIOD_CD",
    "ACTIVE_FLAG_CD",
    "INSPER_EFFECTIVE_DATE",
    "INSPER_EXPIRY_DATE",
    "REPROG_ID_PRE",
    "BUSINESS_TITLE",
    "INSPER_LAYER_NO",
    "INSPER_MAIN_CURRENCY_CD",
    "INSPER_UW_YEAR",
    "INSPER_CONTRACT_ORDER",
    "INSPER_CONTRACT_FORM_CD_SLAV",
    "INSPER_CONTRACT_LODRA_CD_SLAV",
    "INSPER_CONTRACT_COVERAGE_CD_SLAV",
    "INSPER_CLAIM_BASIS_CD",
    "INSPER_LODRA_CD_SLAV",
    "INSPER_LOD_TO_RA_DATE_SLAV",
    "INSPER_COMMENT"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Clear out whatever used to live in columns D:T (nothing did before, but
# make sure stray formatting/values from a wider old layout don't linger)
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(2, $col).Value = $null
}

$ws.Cells.Item(2, 1).Value = 1                  # A2 INSPER_ID_PRE
$ws.Cells.Item(2, 3).Value = "quota_share"       # C2 TYPE_OF_PARTICIPATION_CD
$ws.Cells.Item(2, 5).Value = $true               # E2 ACTIVE_FLAG_CD
$ws.Cells.Item(2, 8).Value = 1                   # H2 REPROG_ID_PRE
$ws.Cells.Item(2, 9).Value = "QS_30"             # I2 BUSINESS_TITLE
$ws.Cells.Item(2, 13).Value = 1                  # M2 INSPER_CONTRACT_ORDER

# ---------------------------------------------------------------------
# 2. "sections" sheet -> rename the first header
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("sections")
$ws3.Cells.Item(1, 1).Value = "BUSINESS_TITLE"
